$wb = $excel.ActiveWorkbook

# Work on the "Repayment schedule" sheet
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new column before column N (14th column), shifting old N/O/P -> O/P/Q
$ws.Columns("N:N").Insert()

# Match the width used for the new column (same as neighboring "Principal" column)
$ws.Columns("N:N").ColumnWidth = 9.83

# Make "Repayment schedule" the active sheet/tab and update its selection
$ws.Activate()
$ws.Range("R5").Select()
